# Logged Week 16 and performed season sim from Week 17
$wb = $excel.ActiveWorkbook

# --- OFF sheet: Home row (row 2) target-depth counts ---
$wsOff = $wb.Worksheets.Item("OFF")
$wsOff.Range("B2").Value = 484
$wsOff.Range("C2").Value = 369
$wsOff.Range("D2").Value = 142
$wsOff.Range("E2").Value = 70

# --- DEF sheet: Home row (row 2) target-depth counts ---
$wsDef = $wb.Worksheets.Item("DEF")
$wsDef.Range("B2").Value = 523
$wsDef.Range("C2").Value = 372
$wsDef.Range("D2").Value = 110
$wsDef.Range("E2").Value = 45
$wsDef.Range("F2").Value = 8
$wsDef.Range("G2").Value = 12
